$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("J1").Value = "depth to water table (m)"

# New column data: depth to water table (m)
$ws.Range("J2").Value = 18
$ws.Range("J3").Value = 8
$ws.Range("J4").Value = 8
$ws.Range("J5").Value = 13.6
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 18
$ws.Range("J8").Value = 8
$ws.Range("J9").Value = 18
$ws.Range("J10").Value = 18
$ws.Range("J11").Value = 13.6
$ws.Range("J12").Value = "NA"
$ws.Range("J13").Value = 13.6
$ws.Range("J14").Value = 13.6

# Column width for new column J
$ws.Columns.Item(10).ColumnWidth = 20

# Sheet view updates (zoom level and active selection cell)
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("I12").Select()
